# Error Calculations and Plots
# This script updates the "missing data" worksheet:
#  - removes the rows for "RM 232" and "SC 92" (they were dropped from the
#    dataset entirely, shifting all subsequent rows up)
#  - updates a number of cells that flip between a present numeric value and
#    a missing value (to reflect the refreshed random "missingness" mask)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped from the sheet.
# Row 26 is "RM 232"; after it is removed, the row that used to be 28
# ("SC 92") becomes row 27.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Apply the remaining cell-level value changes (newly filled-in values and
# newly-blanked/missing values) now that the rows above have shifted up.
$ws.Range("C2").Value = 14.9
$ws.Range("E2").Value = -7.2
$ws.Range("F2").Value = 18.03

$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()

$ws.Range("D4").Value = -15.4
$ws.Range("F4").ClearContents()

$ws.Range("D5").ClearContents()

$ws.Range("C6").ClearContents()

$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()

$ws.Range("C12").Value = 12.5

$ws.Range("E13").Value = -5.3

$ws.Range("C14").ClearContents()

$ws.Range("F18").ClearContents()

$ws.Range("E19").ClearContents()

$ws.Range("C20").Value = 12.5

$ws.Range("C21").Value = 12.7

$ws.Range("C22").ClearContents()
$ws.Range("F22").Value = 16.81

$ws.Range("C23").ClearContents()
$ws.Range("D23").Value = -13.9

$ws.Range("E25").Value = -7.1

$ws.Range("D27").ClearContents()

$ws.Range("E28").Value = -5.9

$ws.Range("D29").Value = -13

$ws.Range("B30").Value = -19.7
$ws.Range("F30").ClearContents()

$ws.Range("C31").Value = 15.3
$ws.Range("E31").ClearContents()

$ws.Range("B32").ClearContents()
$ws.Range("E32").Value = -6.4

$ws.Range("C33").Value = 10.4
